# Appends rows 52-58 to the "Artfynd" sheet (Swedish species-observation
# export) exactly as produced by the upstream automatic-update job.
#
# Columns that hold date/time-looking text (Y, Z, AA, AB) are written with a
# leading apostrophe so Excel stores them as literal text instead of
# auto-converting them to date/time serials - this mirrors how the source
# data is literal text ("2023-09-15", "00:00"), not real Excel dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{
        row = 52
        A = 112111378; B = 82949; C = "Ovaliderad"; D = "NT"; E = 5589
        F = "Rödbrun klubbdyna"; G = "Trichoderma nybergianum"
        H = "(T.Ulvinen & H.L.Chamb.) Jaklitsch & Voglmayr"
        P = "Renkullmyren (Renkullmyren), Jmt"
        Q = 446759.6540099295; R = 7032715.24812395; S = 10
        T = "Jämtland"; U = "Krokom"; V = "Jämtland"; W = "Alsen"
        Y = "2023-09-15"; Z = "00:00"; AA = "2023-09-15"; AB = "00:00"
        AC = $null
        AD = $false; AE = $false; AG = $false
        AW = "Rashid Kadhim"; AX = "Rashid Kadhim"
    },
    @{
        row = 53
        A = 112111398; B = 88966; C = "Ovaliderad"; D = "NT"; E = 5754
        F = "Gultoppig fingersvamp"; G = "Ramaria testaceoflava"
        H = "(Bres.) Corner"
        P = "Renkullmyren (Renkullmyren), Jmt"
        Q = 446739.7436773395; R = 7032704.828598888; S = 10
        T = "Jämtland"; U = "Krokom"; V = "Jämtland"; W = "Alsen"
        Y = "2023-09-15"; Z = "00:00"; AA = "2023-09-15"; AB = "00:00"
        AC = $null
        AD = $false; AE = $false; AG = $false
        AW = "Rashid Kadhim"; AX = "Rashid Kadhim"
    },
    @{
        row = 54
        A = 112110532; B = 88002; C = "Ovaliderad"; D = "VU"; E = 245031
        F = "Borgsjömusseron"; G = "Tricholoma borgsjoeënse"
        H = "Jacobsson & Muskos"
        P = "Svensbergsbäcken (Svensbergsbäcken), Jmt"
        Q = 446764.8121887931; R = 7032863.433365297; S = 10
        T = "Jämtland"; U = "Krokom"; V = "Jämtland"; W = "Alsen"
        Y = "2023-09-15"; Z = "00:00"; AA = "2023-09-15"; AB = "00:00"
        AC = "Längs med en liten stig. På svag sluttning bland kam-och husmossa, revlummer och ekbräken. I närheten finns granvaxskivling, rosa/besk vaxskivling och äggvaxskivling."
        AD = $false; AE = $false; AG = $false
        AW = "Rashid Kadhim"; AX = "Rashid Kadhim"
    },
    @{
        row = 55
        A = 112111386; B = 88946; C = "Ovaliderad"; D = "VU"; E = 256335
        F = "Taggfingersvamp"; G = "Ramaria karstenii"
        H = "(Sacc. & P.Syd.) Corner"
        P = "Renkullmyren (Renkullmyren), Jmt"
        Q = 446733.532634148; R = 7032708.963686833; S = 10
        T = "Jämtland"; U = "Krokom"; V = "Jämtland"; W = "Alsen"
        Y = "2023-09-15"; Z = "00:00"; AA = "2023-09-15"; AB = "00:00"
        AC = $null
        AD = $false; AE = $false; AG = $false
        AW = "Rashid Kadhim"; AX = "Rashid Kadhim"
    },
    @{
        row = 56
        A = 112111388; B = 88956; C = "Ovaliderad"; D = "VU"; E = 5747
        F = "Läderdoftande fingersvamp"; G = "Ramaria safraniolens"
        H = "Christian"
        P = "Renkullmyren (Renkullmyren), Jmt"
        Q = 446733.532634148; R = 7032708.963686833; S = 10
        T = "Jämtland"; U = "Krokom"; V = "Jämtland"; W = "Alsen"
        Y = "2023-09-15"; Z = "00:00"; AA = "2023-09-15"; AB = "00:00"
        AC = $null
        AD = $false; AE = $false; AG = $false
        AW = "Rashid Kadhim"; AX = "Rashid Kadhim"
    },
    @{
        row = 57
        A = 112105381; B = 88899; C = "Ovaliderad"; D = "NT"; E = 3286
        F = "Flattoppad klubbsvamp"; G = "Clavariadelphus truncatus"
        H = "(Quél.) Donk"
        P = "Landverktjärnen (Landverktjärnen), Jmt"
        Q = 446563.6411143647; R = 7032715.829486988; S = 10
        T = "Jämtland"; U = "Krokom"; V = "Jämtland"; W = "Alsen"
        Y = "2023-09-15"; Z = "00:00"; AA = "2023-09-15"; AB = "00:00"
        AC = $null
        AD = $false; AE = $false; AG = $false
        AW = "Rashid Kadhim"; AX = "Rashid Kadhim"
    },
    @{
        row = 58
        A = 112104863; B = 90651; C = "Ovaliderad"; D = "NT"; E = 1968
        F = "Grantaggsvamp"; G = "Bankera violascens"
        H = "(Alb. & Schwein. : Fr.) Pouzar"
        P = "Hökån (Hökån), Jmt"
        Q = 446637.1411376887; R = 7032523.765577726; S = 10
        T = "Jämtland"; U = "Krokom"; V = "Jämtland"; W = "Alsen"
        Y = "2023-09-15"; Z = "00:00"; AA = "2023-09-15"; AB = "00:00"
        AC = $null
        AD = $false; AE = $false; AG = $false
        AW = "Rashid Kadhim"; AX = "Rashid Kadhim"
    }
)

# Columns holding literal date/time text - must be forced to text so Excel
# doesn't reinterpret e.g. "2023-09-15" as a date serial number.
$dateTextCols = @("Y", "Z", "AA", "AB")

foreach ($r in $rows) {
    $rowNum = $r.row

    $ws.Range("A$rowNum").Value = $r.A
    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("D$rowNum").Value = $r.D
    $ws.Range("E$rowNum").Value = $r.E
    $ws.Range("F$rowNum").Value = $r.F
    $ws.Range("G$rowNum").Value = $r.G
    $ws.Range("H$rowNum").Value = $r.H

    $ws.Range("P$rowNum").Value = $r.P
    $ws.Range("Q$rowNum").Value = $r.Q
    $ws.Range("R$rowNum").Value = $r.R
    $ws.Range("S$rowNum").Value = $r.S
    $ws.Range("T$rowNum").Value = $r.T
    $ws.Range("U$rowNum").Value = $r.U
    $ws.Range("V$rowNum").Value = $r.V
    $ws.Range("W$rowNum").Value = $r.W

    foreach ($col in $dateTextCols) {
        $cell = $ws.Range("$col$rowNum")
        $cell.NumberFormat = "@"
        $cell.Value = $r[$col]
    }

    if ($r.AC) {
        $ws.Range("AC$rowNum").Value = $r.AC
    }

    $ws.Range("AD$rowNum").Value = $r.AD
    $ws.Range("AE$rowNum").Value = $r.AE
    $ws.Range("AG$rowNum").Value = $r.AG

    $ws.Range("AW$rowNum").Value = $r.AW
    $ws.Range("AX$rowNum").Value = $r.AX
}
